# The workbook contains a single sheet of daily Mango price records
# (Vega Central Mapocho de Santiago) ordered by row. This commit adds a
# new daily record, inserted as row 354, which pushes all the existing
# records that used to occupy rows 354-449 down by one row (to 355-450).
#
# Net effect on the sheet dimension: A1:T449 -> A1:T450.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 354; Excel shifts rows 354..449 down to 355..450
# and the surrounding formatting (e.g. the date-style column D) carries over.
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row 354 with the new record's data.
$ws.Cells.Item(354, 1).Value  = 9
$ws.Cells.Item(354, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(354, 3).Value  = "Metropolitana"
$ws.Cells.Item(354, 4).Value  = 44754
$ws.Cells.Item(354, 5).Value  = 13
$ws.Cells.Item(354, 6).Value  = "Fruta"
$ws.Cells.Item(354, 7).Value  = 100108
$ws.Cells.Item(354, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(354, 9).Value  = 100108002
$ws.Cells.Item(354, 10).Value = "Mango"
$ws.Cells.Item(354, 11).Value = "Sin especificar"
$ws.Cells.Item(354, 12).Value = "Primera"
$ws.Cells.Item(354, 13).Value = 500
$ws.Cells.Item(354, 14).Value = 7000
$ws.Cells.Item(354, 15).Value = 7500
$ws.Cells.Item(354, 16).Value = 7280
$ws.Cells.Item(354, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(354, 18).Value = "Brasil"
$ws.Cells.Item(354, 19).Value = 1820
$ws.Cells.Item(354, 20).Value = 4
